$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 181 (shifts existing rows 181..225 down to 182..226)
$ws.Rows.Item(181).Insert()

# Populate the new row 181 with the new price-record data
$ws.Cells.Item(181, 1).Value = 7
$ws.Cells.Item(181, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(181, 3).Value = "Ñuble"
$ws.Cells.Item(181, 4).Value = 44943
$ws.Cells.Item(181, 5).Value = 16
$ws.Cells.Item(181, 6).Value = 100112045
$ws.Cells.Item(181, 7).Value = "Zapallo"
$ws.Cells.Item(181, 8).Value = "Camote"
$ws.Cells.Item(181, 9).Value = "1a (cosecha)"
$ws.Cells.Item(181, 10).Value = 300
$ws.Cells.Item(181, 11).Value = 600
$ws.Cells.Item(181, 12).Value = 650
$ws.Cells.Item(181, 13).Value = 625
$ws.Cells.Item(181, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(181, 15).Value = "Región del Maule"
$ws.Cells.Item(181, 16).Value = 625
$ws.Cells.Item(181, 17).Value = 1
$ws.Cells.Item(181, 18).Value = "Hortaliza"
